$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Coin name / Link swaps (rows 33/34 and 42/43) ---
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# --- Update Price (column D) values as text, preserving default style ---
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "26.849.96"
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.801.59"
$r.Style = "Normal"
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.001"
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "309.15"
$r.Style = "Normal"
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.4656"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.3691"
$r.Style = "Normal"
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.07365"
$r.Style = "Normal"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.8696"
$r.Style = "Normal"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "20.34"
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.778.52"
$r.Style = "Normal"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "5.354"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "92.41"
$r.Style = "Normal"
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "6.497"
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "1.001"
$r.Style = "Normal"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.000008690"
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "26.852.87"
$r.Style = "Normal"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "5.286"
$r.Style = "Normal"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "10.60"
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "2.004.27"
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "1.908"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "151.67"
$r.Style = "Normal"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "18.31"
$r.Style = "Normal"
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.128"
$r.Style = "Normal"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "5.252"
$r.Style = "Normal"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "116.11"
$r.Style = "Normal"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.08913"
$r.Style = "Normal"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.7588"
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.149"
$r.Style = "Normal"
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "2.929"
$r.Style = "Normal"
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "4.459"
$r.Style = "Normal"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.01951"
$r.Style = "Normal"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.05247"
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "2.926"
$r.Style = "Normal"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "7.215"
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.5292"
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "2.359"
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.1659"
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "8.483"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.4999"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "10.34"
$r.Style = "Normal"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "104.20"
$r.Style = "Normal"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.662"
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.06283"
$r.Style = "Normal"

# --- Update Volume(1h) (column E) percentage text values ---
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +3.94%  "
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("E23").Value = "  -3.59%  "
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("E28").Value = "  -8.62%  "
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("E33").Value = "  -4.54%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("E51").Value = "  -2.00%  "
